$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 30)
$ws.Range("D2").Value = [double]"0.9999998079147523"
$ws.Range("E2").Value = [double]"0.9999998079147523"

# Row 3 (Control 11)
$ws.Range("D3").Value = [double]"0.9999999977843042"
$ws.Range("E3").Value = [double]"0.9999999977843042"

# Row 4 (Control 3)
$ws.Range("D4").Value = [double]"1.236224896776247E-27"
$ws.Range("E4").Value = [double]"1.236224896776247E-27"

# Row 5 (Control 38)
$ws.Range("D5").Value = [double]"0.9901112600426142"
$ws.Range("E5").Value = [double]"0.9901112600426142"

# Row 6 (Control 29)
$ws.Range("D6").Value = [double]"0.9999919240617908"
$ws.Range("E6").Value = [double]"0.9999919240617908"

# Row 7 (MDD 41)
$ws.Range("D7").Value = [double]"0.002131138354101981"
$ws.Range("E7").Value = [double]"0.9978688616458981"

# Row 8 (MDD 8)
$ws.Range("D8").Value = [double]"0.9999999999998541"
$ws.Range("E8").Value = [double]"1.458833054357456E-13"

# Row 9 (MDD 15)
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = [double]"0.9999999999964497"
$ws.Range("E9").Value = [double]"3.550271188146326E-12"

# Row 11 (MDD 33)
$ws.Range("F11").Value = [double]"5.788710594177246"
$ws.Range("G11").Value = [double]"0.5"
